$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (column D) values are written as text, not auto-converted to numbers,
# then restore the default "Normal" style so no stray formatting is introduced.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.460.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.504.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.10%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.72'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('E10').Value = '  +4.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.400'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.104.40'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.69'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.50%  '
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.498.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.23%  '
$ws.Range('E16').Value = '  +4.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.518.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.29'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '395.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.568'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '75.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('E25').Value = '  +9.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.642.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.191'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.995'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.38%  '
$ws.Range('E32').Value = '  +6.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.70%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.66%  '
$ws.Range('B36').Value = 'EnergySwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.66'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +30.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '173.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.539.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.14%  '
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('E42').Value = '  +4.25%  '
$ws.Range('E43').Value = '  +8.21%  '
$ws.Range('E44').Value = '  +4.89%  '
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('E46').Value = '  +10.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.603.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.33%  '
$ws.Range('E49').Value = '  +12.56%  '
$ws.Range('E50').Value = '  +3.04%  '
$ws.Range('E51').Value = '  +5.49%  '
